# Add a new wave of data (column for the "16. 3. 2021" survey date) to both
# sheets, right after the existing last date column, and bump the
# "aktualizace" date mentioned in the two summary-title cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "data": new column AA (right after Z) = "16. 3. 2021"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Copy the header formatting from the last existing date column (Z1) onto
# the new one (AA1), then overwrite its value with the new date.
$ws1.Range("Z1").Copy($ws1.Range("AA1"))
$ws1.Cells.Item(1, 27).Value = "16. 3. 2021"

$aaValues = @(0.4,0.24,0.13,0.12,0.11,0.3,0.22,0.24,0.09,0.15,0.3,0.22,0.13,0.18,0.17,0.53,0.26,0.08,0.08,0.05,0.3,0.16,0.27,0.09,0.18,0.41,0.2,0.14,0.11,0.14,0.42,0.28,0.09,0.13,0.08,0.3,0.22,0.15,0.18,0.15,0.42,0.27,0.11,0.1,0.1,0.48,0.2,0.15,0.08,0.09,0.34,0.21,0.15,0.14,0.16,0.45,0.26,0.12,0.1,0.07000000000000001,0.42,0.21,0.14,0.12,0.11,0.37,0.29,0.11,0.12,0.11,0.35,0.25,0.13,0.14,0.13,0.39,0.25,0.17,0.09,0.1,0.27,0.17,0.28,0.1,0.18,0.22,0.2,0.15,0.21,0.22,0.51,0.26,0.07000000000000001,0.09,0.07000000000000001,0.33,0.27,0.21,0.08,0.11,0.39,0.24,0.11,0.14,0.12,0.55,0.27,0.08,0.08,0.02)

for ($i = 0; $i -lt $aaValues.Length; $i++) {
    $ws1.Cells.Item($i + 2, 27).Value = $aaValues[$i]
}

# Update the "% respondentu" title (row 112, col A): bump the trailing date
$titleCell1 = $ws1.Cells.Item(112, 1)
$titleCell1.Value = ($titleCell1.Value2 -replace "9\. 3\. 2021$", "23. 3. 2021")

# ---------------------------------------------------------------
# Sheet "pocetR": new column Z (right after Y) = "16. 3. 2021"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Copy the header formatting from the last existing date column (Y1) onto
# the new one (Z1), then overwrite its value with the new date.
$ws2.Range("Y1").Copy($ws2.Range("Z1"))
$ws2.Cells.Item(1, 26).Value = "16. 3. 2021"

$zValues = @(1872,431,666,775,316,619,937,627,692,553,909,963,978,425,225,244,220,346,344,211,321,430)

for ($i = 0; $i -lt $zValues.Length; $i++) {
    $ws2.Cells.Item($i + 2, 26).Value = $zValues[$i]
}

# Blank trailing footer cell under the new column, matching the rest of row 24
$ws2.Range("Y24").Copy($ws2.Range("Z24"))

# Update the "velikost dotazaneho souboru" title (row 24, col A): bump the trailing date
$titleCell2 = $ws2.Cells.Item(24, 1)
$titleCell2.Value = ($titleCell2.Value2 -replace "9\. 3\. 2021$", "23. 3. 2021")

Write-Output "edit complete"
